{"js": "// Remove author from metadata.\n// The document's built-in \"Author\" property (docProps/core.xml <dc:creator>)\n// is cleared, matching the commit \"Remove author from metadata.\"\nconst props = context.document.properties;\nprops.load(\"author\");\nawait context.sync();\n\nprops.author = \"\";\nawait context.sync();\n", "ps1": "# Remove author from metadata.\n# Clears the document's built-in Author property (docProps/core.xml <dc:creator>),\n# matching the commit \"Remove author from metadata.\"\n$d = $word.ActiveDocument\n$d.Author = \"\"\n"}
